$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 200001180
$ws.Range("I18").Value = 299.66666
$ws.Range("J18").Value = 500002500
$ws.Range("K18").Value = 299.66666
$ws.Range("L18").Value = 500002500
$ws.Range("M18").Value = -15.66665999999998
$ws.Range("N18").Value = -500003068

$ws.Range("H32").Value = 41666816
$ws.Range("I32").Value = 41666816
$ws.Range("K32").Value = 41666816
$ws.Range("M32").Value = -41666490

$ws.Range("H40").Value = 4069.6428
$ws.Range("I40").Value = 10000
$ws.Range("J40").Value = 3613.4614
$ws.Range("K40").Value = 10000
$ws.Range("L40").Value = 3613.4614
$ws.Range("M40").Value = -9825
$ws.Range("N40").Value = -3963.4614

$ws.Range("H64").Value = 142861520
$ws.Range("J64").Value = 250004290
$ws.Range("L64").Value = 250004290
$ws.Range("N64").Value = -250004786

$ws.Range("H67").Value = 142861520
$ws.Range("J67").Value = 250004290
$ws.Range("L67").Value = 250004290
$ws.Range("N67").Value = -250006006

$ws.Range("H132").Value = 10020.2
$ws.Range("I132").Value = 10628.786
$ws.Range("K132").Value = 31886.358
$ws.Range("M132").Value = -29356.358

$ws.Range("H137").Value = 912134
$ws.Range("I137").Value = 1391042.1
$ws.Range("J137").Value = 4729.2104
$ws.Range("K137").Value = 4173126.3
$ws.Range("L137").Value = 14187.6312
$ws.Range("M137").Value = -4170576.3
$ws.Range("N137").Value = -19287.6312

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1398.0667
$ws.Range("I2").Value = 1357
$ws.Range("K2").Value = 1357
$ws.Range("M2").Value = -1244

$ws.Range("H32").Value = 1728785.4
$ws.Range("I32").Value = 772738.25
$ws.Range("K32").Value = 772738.25
$ws.Range("M32").Value = -772451.25

$ws.Range("H61").Value = 6253.091
$ws.Range("I61").Value = 2697.75
$ws.Range("K61").Value = 2697.75
$ws.Range("M61").Value = -2485.75

$ws.Range("H63").Value = 1498.5
$ws.Range("I63").Value = 1498.5
$ws.Range("K63").Value = 1498.5
$ws.Range("M63").Value = -812.5

$ws.Range("H66").Value = 1498.5
$ws.Range("I66").Value = 1498.5
$ws.Range("K66").Value = 7492.5
$ws.Range("M66").Value = -4060.5

$ws.Range("H74").Value = 2089.9119
$ws.Range("I74").Value = 1603.1428
$ws.Range("K74").Value = 1603.1428
$ws.Range("M74").Value = -729.1428000000001

$ws.Range("H77").Value = 2089.9119
$ws.Range("I77").Value = 1603.1428
$ws.Range("K77").Value = 8015.714
$ws.Range("M77").Value = -3647.714

$ws.Range("H102").Value = 2686.8948
$ws.Range("I102").Value = 2075.2856
$ws.Range("K102").Value = 2075.2856
$ws.Range("M102").Value = -453.2856000000002

$ws.Range("H116").Value = 1398.0667
$ws.Range("I116").Value = 1357
$ws.Range("K116").Value = 1357
$ws.Range("M116").Value = 937

$ws.Range("H122").Value = 1801.0625
$ws.Range("I122").Value = 1787.3334
$ws.Range("K122").Value = 5362.0002
$ws.Range("M122").Value = -2912.0002

$ws.Range("H132").Value = 2957.3215
$ws.Range("I132").Value = 1950.579
$ws.Range("J132").Value = 5082.6665
$ws.Range("K132").Value = 5851.737
$ws.Range("L132").Value = 15247.9995
$ws.Range("M132").Value = -3321.737
$ws.Range("N132").Value = -20307.9995

$ws.Range("H136").Value = 6253.091
$ws.Range("I136").Value = 2697.75
$ws.Range("K136").Value = 8093.25
$ws.Range("M136").Value = -5543.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1398.0667
$ws.Range("I3").Value = 1357
$ws.Range("K3").Value = 1357
$ws.Range("M3").Value = -1243

$ws.Range("H105").Value = 15295806
$ws.Range("I105").Value = 1430254.2
$ws.Range("K105").Value = 1430254.2
$ws.Range("M105").Value = -1428507.2

$ws.Range("H134").Value = 1412.5
$ws.Range("I134").Value = 1421.8334
$ws.Range("K134").Value = 4265.5002
$ws.Range("M134").Value = -1730.5002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 195.57143
$ws.Range("I7").Value = 122.5
$ws.Range("K7").Value = 122.5
$ws.Range("M7").Value = -9.5

$ws.Range("H31").Value = 4171777.8
$ws.Range("I31").Value = 2892.5557
$ws.Range("J31").Value = 5958443
$ws.Range("K31").Value = 2892.5557
$ws.Range("L31").Value = 5958443
$ws.Range("M31").Value = -2597.5557
$ws.Range("N31").Value = -5959033

$ws.Range("H34").Value = 4171777.8
$ws.Range("I34").Value = 2892.5557
$ws.Range("J34").Value = 5958443
$ws.Range("K34").Value = 2892.5557
$ws.Range("L34").Value = 5958443
$ws.Range("M34").Value = -2690.5557
$ws.Range("N34").Value = -5958847

$ws.Range("H99").Value = 3904
$ws.Range("I99").Value = 2868.125
$ws.Range("J99").Value = 6666.3335
$ws.Range("K99").Value = 2868.125
$ws.Range("L99").Value = 6666.3335
$ws.Range("M99").Value = -1370.125
$ws.Range("N99").Value = -9662.333500000001

$ws.Range("H107").Value = 2273391.2
$ws.Range("I107").Value = 3125308.5
$ws.Range("J107").Value = 1611.5
$ws.Range("K107").Value = 3125308.5
$ws.Range("L107").Value = 1611.5
$ws.Range("M107").Value = -3123388.5
$ws.Range("N107").Value = -5451.5

$ws.Range("H126").Value = 3904
$ws.Range("I126").Value = 2868.125
$ws.Range("J126").Value = 6666.3335
$ws.Range("K126").Value = 8604.375
$ws.Range("L126").Value = 19999.0005
$ws.Range("M126").Value = -6134.375
$ws.Range("N126").Value = -24939.0005

$ws.Range("H132").Value = 3707.22
$ws.Range("I132").Value = 3370.111
$ws.Range("J132").Value = 4574.0713
$ws.Range("K132").Value = 10110.333
$ws.Range("L132").Value = 13722.2139
$ws.Range("M132").Value = -7580.332999999999
$ws.Range("N132").Value = -18782.2139

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 4913.933
$ws.Range("I139").Value = 1851.8
$ws.Range("J139").Value = 6445
$ws.Range("K139").Value = 5555.4
$ws.Range("L139").Value = 19335
$ws.Range("M139").Value = -415.3999999999996
$ws.Range("N139").Value = -29615

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 90913370
$ws.Range("I80").Value = 250002140
$ws.Range("J80").Value = 5496.857
$ws.Range("K80").Value = 250002140
$ws.Range("L80").Value = 5496.857
$ws.Range("M80").Value = -250001142
$ws.Range("N80").Value = -7492.857

$ws.Range("H83").Value = 90913370
$ws.Range("I83").Value = 250002140
$ws.Range("J83").Value = 5496.857
$ws.Range("K83").Value = 1250010700
$ws.Range("L83").Value = 27484.285
$ws.Range("M83").Value = -1250005708
$ws.Range("N83").Value = -37468.285

$ws.Range("H97").Value = 3570
$ws.Range("I97").Value = 3587
$ws.Range("J97").Value = 3400
$ws.Range("K97").Value = 3587
$ws.Range("L97").Value = 3400
$ws.Range("M97").Value = -3091
$ws.Range("N97").Value = -4392

$ws.Range("H102").Value = 7772.1562
$ws.Range("I102").Value = 1192.5714
$ws.Range("K102").Value = 1192.5714
$ws.Range("M102").Value = 429.4286

$ws.Range("H132").Value = 2349.9644
$ws.Range("I132").Value = 2190.7778
$ws.Range("K132").Value = 6572.3334
$ws.Range("M132").Value = -4042.3334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2179.8572
$ws.Range("I7").Value = 1647.3334
$ws.Range("J7").Value = 3138.4
$ws.Range("K7").Value = 1647.3334
$ws.Range("L7").Value = 3138.4
$ws.Range("M7").Value = -1535.3334
$ws.Range("N7").Value = -3362.4

$ws.Range("H22").Value = 89286616
$ws.Range("I22").Value = 8929560
$ws.Range("J22").Value = 250000720
$ws.Range("K22").Value = 8929560
$ws.Range("L22").Value = 250000720
$ws.Range("M22").Value = -8929265
$ws.Range("N22").Value = -250001310

$ws.Range("H27").Value = 89286616
$ws.Range("I27").Value = 8929560
$ws.Range("J27").Value = 250000720
$ws.Range("K27").Value = 8929560
$ws.Range("L27").Value = 250000720
$ws.Range("M27").Value = -8929453
$ws.Range("N27").Value = -250000934

$ws.Range("H46").Value = 3006.3845
$ws.Range("J46").Value = 5664.6665
$ws.Range("L46").Value = 5664.6665
$ws.Range("N46").Value = -6040.6665

$ws.Range("H61").Value = 899.6
$ws.Range("I61").Value = 899.75
$ws.Range("K61").Value = 899.75
$ws.Range("M61").Value = -697.75

$ws.Range("H113").Value = 899.6
$ws.Range("I113").Value = 899.75
$ws.Range("K113").Value = 899.75
$ws.Range("M113").Value = 1270.25

$ws.Range("H122").Value = 5913.294
$ws.Range("J122").Value = 8128
$ws.Range("L122").Value = 24384
$ws.Range("N122").Value = -29284

$ws.Range("H126").Value = 2179.8572
$ws.Range("I126").Value = 1647.3334
$ws.Range("J126").Value = 3138.4
$ws.Range("K126").Value = 4942.0002
$ws.Range("L126").Value = 9415.200000000001
$ws.Range("M126").Value = -2472.0002
$ws.Range("N126").Value = -14355.2

$ws.Range("H136").Value = 2726.8845
$ws.Range("I136").Value = 2057.45
$ws.Range("K136").Value = 6172.349999999999
$ws.Range("M136").Value = -3622.349999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 16667762
$ws.Range("I122").Value = 1280.2727
$ws.Range("K122").Value = 3840.8181
$ws.Range("M122").Value = -1390.8181

$ws.Range("H126").Value = 1902.8462
$ws.Range("I126").Value = 1903.4546
$ws.Range("K126").Value = 5710.3638
$ws.Range("M126").Value = -3240.3638

$ws.Range("H132").Value = 2051.1333
$ws.Range("I132").Value = 1797.4546
$ws.Range("J132").Value = 2748.75
$ws.Range("K132").Value = 5392.3638
$ws.Range("L132").Value = 8246.25
$ws.Range("M132").Value = -2862.3638
$ws.Range("N132").Value = -13306.25

$ws.Range("H136").Value = 2589.3845
$ws.Range("I136").Value = 2660.96
$ws.Range("K136").Value = 7982.88
$ws.Range("M136").Value = -5432.88
